# Update the dSF column (F) values for a specific set of rows,
# matching the "repull data, push all data, mean calculation" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    5  = -1
    11 = -10
    15 = -3
    16 = -1
    27 = 1
    31 = -1
    34 = -1
    37 = 4
    40 = 0
    49 = -2
    55 = -3
    62 = -2
    63 = -4
    65 = -3
    67 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
